# Simulated Wild Card round and logged it
$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$rushing = $wb.Worksheets.Item("Rushing")

# Row 2: D.Lock
$rushing.Range("D2").Value = 4
$rushing.Range("E2").Value = 6
$rushing.Range("F2").Value = 2

# Row 4: M.Gordon
$rushing.Range("C4").Value = 110
$rushing.Range("D4").Value = 77
$rushing.Range("F4").Value = 37

# Row 5: J.Williams
$rushing.Range("C5").Value = 112
$rushing.Range("D5").Value = 67
$rushing.Range("E5").Value = 23
$rushing.Range("F5").Value = 29

# --- Receiving sheet ---
$receiving = $wb.Worksheets.Item("Receiving")

# Row 2: M.Gordon
$receiving.Range("C2").Value = 36
$receiving.Range("D2").Value = 25

# Row 3: J.Williams
$receiving.Range("C3").Value = 51
$receiving.Range("D3").Value = 41

# Row 5: C.Sutton
$receiving.Range("C5").Value = 65
$receiving.Range("D5").Value = 49

# Row 6: J.Jeudy
$receiving.Range("C6").Value = 76
$receiving.Range("D6").Value = 60
$receiving.Range("E6").Value = 27
$receiving.Range("F6").Value = 18

# Row 7: T.Patrick
$receiving.Range("C7").Value = 64
$receiving.Range("D7").Value = 43
$receiving.Range("E7").Value = 21
$receiving.Range("F7").Value = 10
$receiving.Range("G7").Value = 11

# Row 11: N.Fant
$receiving.Range("C11").Value = 82
$receiving.Range("D11").Value = 65
$receiving.Range("G11").Value = 16
